$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$values = @{
    2  = "Not worse"
    3  = "A little worse"
    4  = "A little worse"
    5  = "A little worse"
    6  = "A little worse"
    7  = "A lot worse"
    8  = "Not worse"
    9  = "Somewhat worse"
    10 = "Somewhat worse"
    11 = "A little worse"
    12 = "A little worse"
    13 = "Not worse"
    14 = "Not worse"
    15 = "A lot worse"
}

foreach ($row in $values.Keys) {
    $ws.Range("C$row").Value = $values[$row]
}

$ws.Range("C16").Select()
